$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 52; this shifts all existing rows (52..110) down to (53..111)
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with the new record's data
$ws.Cells.Item(52, 1).Value = 4
$ws.Cells.Item(52, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(52, 3).Value = "Los Lagos"
$ws.Cells.Item(52, 4).Value = 44413
$ws.Cells.Item(52, 5).Value = 10
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100108
$ws.Cells.Item(52, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(52, 9).Value = 100108005
$ws.Cells.Item(52, 10).Value = "Piña"
$ws.Cells.Item(52, 11).Value = "Caramelo"
$ws.Cells.Item(52, 12).Value = "Segunda"
$ws.Cells.Item(52, 13).Value = 120
$ws.Cells.Item(52, 14).Value = 21000
$ws.Cells.Item(52, 15).Value = 21000
$ws.Cells.Item(52, 16).Value = 21000
$ws.Cells.Item(52, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(52, 18).Value = "Ecuador"
$ws.Cells.Item(52, 19).Value = 1500
$ws.Cells.Item(52, 20).Value = 14
